# Applies the "fixing bug laporan, pembayaran ujian, jenis ujian, revisi
# database" revision to the student-fee report workbook.
#
# Sheet "Data 2"  -> student is now in Kelas 12 (XII_ak_A): SPP figures for
#                    Kelas 10/11 raised, Kelas 12 section filled in with real
#                    figures, and three new summary rows (kelas 10/11/12 +
#                    grand total) appended.
# Sheet "Data 1"  -> student is in Kelas 11 (no Kelas 12 section anymore):
#                    SPP figures for Kelas 10/11 corrected and the trailing
#                    Kelas 12 block is replaced by a compact totals block.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Data 2"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Data 2")

# Student moved up to Kelas XII
$ws1.Range("B8").Value = "XII_ak_A"

# DPP (still outstanding) shown as "0,-" instead of a paid amount
$ws1.Range("B9").Value = "0,-"

# Kelas 10 charges revised upward
$ws1.Range("B14").Value = 300000
$ws1.Range("B15").Value = 150000
$ws1.Range("B16").Value = 300000
$ws1.Range("B17").Value = 150000
$ws1.Range("B18").Value = 1260000

# Kelas 11 charges revised upward
$ws1.Range("B23").Value = 300000
$ws1.Range("B24").Value = 150000
$ws1.Range("B25").Value = 300000
$ws1.Range("B26").Value = 150000
$ws1.Range("B27").Value = 1260000

# Kelas 12 block now filled in with the real class + figures
$ws1.Range("B30").Value = "XII_ak_A"
$ws1.Range("B31").Value = 360000
$ws1.Range("B33").Value = 150000
$ws1.Range("B35").Value = 150000
$ws1.Range("B36").Value = 625000
$ws1.Range("B37").Value = 1285000

# Overall DPP still outstanding
$ws1.Range("B40").Value = "0,-"

# Row 41 becomes "Total Tanggungan kelas 10" (was "total Tanggugan Biaya")
$ws1.Range("A41").Value = "Total Tanggungan kelas 10"
$ws1.Range("B41").Value = 1260000

# New rows 42-44: per-class totals + grand total, styled like row 41
$ws1.Range("A41:B41").Copy() | Out-Null
$ws1.Range("A42:B44").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws1.Range("A42").Value = "Total Tanggungan kelas 11"
$ws1.Range("B42").Value = 1260000
$ws1.Range("A43").Value = "Total Tanggungan kelas 12"
$ws1.Range("B43").Value = 1285000
$ws1.Range("A44").Value = "total Tanggugan Biaya"
$ws1.Range("B44").Value = 3805000

# Keep the view's selection in sync with the new used range
[void]$ws1.Range("B3:B44").Select()

# ---------------------------------------------------------------------
# Sheet "Data 1"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Data 1")

# Kelas 10 charges revised
$ws2.Range("B14").Value = 300000
$ws2.Range("B16").Value = 300000
$ws2.Range("B17").Value = "0,-"
$ws2.Range("B18").Value = 760000

# Kelas 11 charges revised
$ws2.Range("B23").Value = 300000
$ws2.Range("B24").Value = 150000
$ws2.Range("B25").Value = 300000
$ws2.Range("B26").Value = 150000
$ws2.Range("B27").Value = 1140000

# The trailing Kelas 12 block (old rows 29-41) is replaced by a compact
# totals block (new rows 29-34). Pull the formats from the existing
# "blank separator" / "label" / "label + amount" template rows before the
# old rows are removed.
$ws2.Range("A28:B28").Copy() | Out-Null
$ws2.Range("A29:B29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws2.Range("A39:B39").Copy() | Out-Null
$ws2.Range("A30:B30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws2.Range("A40:B40").Copy() | Out-Null
$ws2.Range("A31:B31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws2.Range("A41:B41").Copy() | Out-Null
$ws2.Range("A32:B34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 29 becomes a blank separator row (values cleared)
$ws2.Range("A29").Value = $null
$ws2.Range("B29").Value = $null

# Row 30: "Total Tanggungan" label
$ws2.Range("A30").Value = "Total Tanggungan"

# Row 31: outstanding DPP
$ws2.Range("A31").Value = "DPP"
$ws2.Range("B31").Value = 200000

# Row 32-34: per-class totals + grand total
$ws2.Range("A32").Value = "Total Tanggungan kelas 10"
$ws2.Range("B32").Value = 760000
$ws2.Range("A33").Value = "Total Tanggungan kelas 11"
$ws2.Range("B33").Value = 1140000
$ws2.Range("A34").Value = "total Tanggugan Biaya"
$ws2.Range("B34").Value = 2100000

# Drop the now-obsolete Kelas 12 rows (old 35-41), seven rows in total
$ws2.Range("A35:A41").EntireRow.Delete() | Out-Null

# The Kelas 12 heading used to live in a merged A29:B29 cell; that heading
# is gone now that row 29 is a blank separator, so un-merge it.
$ws2.Range("A29:B29").UnMerge() | Out-Null

# Keep the view's selection in sync with the new used range
[void]$ws2.Range("B3:B34").Select()
